$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(82, 8).Value = 1105.3334
$ws.Cells.Item(82, 9).Value = 1105.3334
$ws.Cells.Item(82, 11).Value = 3316.0002
$ws.Cells.Item(82, 13).Value = -2910.0002
$ws.Cells.Item(85, 8).Value = 1105.3334
$ws.Cells.Item(85, 9).Value = 1105.3334
$ws.Cells.Item(85, 11).Value = 3316.0002
$ws.Cells.Item(85, 13).Value = -1912.0002
$ws.Cells.Item(111, 8).Value = 504500.5
$ws.Cells.Item(111, 9).Value = 9000
$ws.Cells.Item(111, 11).Value = 27000
$ws.Cells.Item(111, 13).Value = -23933
$ws.Cells.Item(113, 8).Value = 6669556.5
$ws.Cells.Item(113, 9).Value = 9526395
$ws.Cells.Item(113, 10).Value = 3600
$ws.Cells.Item(113, 11).Value = 9526395
$ws.Cells.Item(113, 12).Value = 3600
$ws.Cells.Item(113, 13).Value = -9523141
$ws.Cells.Item(113, 14).Value = -10108
$ws.Cells.Item(137, 8).Value = 1723.5682
$ws.Cells.Item(137, 9).Value = 1661.5555
$ws.Cells.Item(137, 10).Value = 1822.0588
$ws.Cells.Item(137, 11).Value = 4984.666499999999
$ws.Cells.Item(137, 12).Value = 5466.1764
$ws.Cells.Item(137, 13).Value = -2434.666499999999
$ws.Cells.Item(137, 14).Value = -10566.1764
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2584.9092
$ws.Cells.Item(2, 9).Value = 2224.2
$ws.Cells.Item(2, 10).Value = 2885.5
$ws.Cells.Item(2, 11).Value = 2224.2
$ws.Cells.Item(2, 12).Value = 2885.5
$ws.Cells.Item(2, 13).Value = -2111.2
$ws.Cells.Item(2, 14).Value = -3111.5
$ws.Cells.Item(32, 8).Value = 7215.7705
$ws.Cells.Item(32, 9).Value = 5243.7407
$ws.Cells.Item(32, 11).Value = 5243.7407
$ws.Cells.Item(32, 13).Value = -4956.7407
$ws.Cells.Item(45, 8).Value = 1188.9
$ws.Cells.Item(45, 9).Value = 1048.625
$ws.Cells.Item(45, 11).Value = 1048.625
$ws.Cells.Item(45, 13).Value = -671.625
$ws.Cells.Item(61, 8).Value = 7261.5713
$ws.Cells.Item(61, 9).Value = 8935.214
$ws.Cells.Item(61, 10).Value = 3914.2856
$ws.Cells.Item(61, 11).Value = 8935.214
$ws.Cells.Item(61, 12).Value = 3914.2856
$ws.Cells.Item(61, 13).Value = -8723.214
$ws.Cells.Item(61, 14).Value = -4338.2856
$ws.Cells.Item(116, 8).Value = 2584.9092
$ws.Cells.Item(116, 9).Value = 2224.2
$ws.Cells.Item(116, 10).Value = 2885.5
$ws.Cells.Item(116, 11).Value = 2224.2
$ws.Cells.Item(116, 12).Value = 2885.5
$ws.Cells.Item(116, 13).Value = 69.80000000000018
$ws.Cells.Item(116, 14).Value = -7473.5
$ws.Cells.Item(122, 8).Value = 694339.75
$ws.Cells.Item(122, 9).Value = 802599.4
$ws.Cells.Item(122, 10).Value = 1478
$ws.Cells.Item(122, 11).Value = 2407798.2
$ws.Cells.Item(122, 12).Value = 4434
$ws.Cells.Item(122, 13).Value = -2405348.2
$ws.Cells.Item(122, 14).Value = -9334
$ws.Cells.Item(132, 8).Value = 3707239
$ws.Cells.Item(132, 9).Value = 2142.5293
$ws.Cells.Item(132, 10).Value = 10005903
$ws.Cells.Item(132, 11).Value = 6427.5879
$ws.Cells.Item(132, 12).Value = 30017709
$ws.Cells.Item(132, 13).Value = -3897.5879
$ws.Cells.Item(132, 14).Value = -30022769
$ws.Cells.Item(136, 8).Value = 7261.5713
$ws.Cells.Item(136, 9).Value = 8935.214
$ws.Cells.Item(136, 10).Value = 3914.2856
$ws.Cells.Item(136, 11).Value = 26805.642
$ws.Cells.Item(136, 12).Value = 11742.8568
$ws.Cells.Item(136, 13).Value = -24255.642
$ws.Cells.Item(136, 14).Value = -16842.8568
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2584.9092
$ws.Cells.Item(3, 9).Value = 2224.2
$ws.Cells.Item(3, 10).Value = 2885.5
$ws.Cells.Item(3, 11).Value = 2224.2
$ws.Cells.Item(3, 12).Value = 2885.5
$ws.Cells.Item(3, 13).Value = -2110.2
$ws.Cells.Item(3, 14).Value = -3113.5
$ws.Cells.Item(107, 8).Value = 1146.2
$ws.Cells.Item(107, 9).Value = 1017.6923
$ws.Cells.Item(107, 11).Value = 1017.6923
$ws.Cells.Item(107, 13).Value = 902.3077
$ws.Cells.Item(134, 8).Value = 4292.275
$ws.Cells.Item(134, 9).Value = 4566.394
$ws.Cells.Item(134, 10).Value = 3000
$ws.Cells.Item(134, 11).Value = 13699.182
$ws.Cells.Item(134, 12).Value = 9000
$ws.Cells.Item(134, 13).Value = -11164.182
$ws.Cells.Item(134, 14).Value = -14070
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2742.5
$ws.Cells.Item(132, 9).Value = 1500
$ws.Cells.Item(132, 10).Value = 2991
$ws.Cells.Item(132, 11).Value = 4500
$ws.Cells.Item(132, 12).Value = 8973
$ws.Cells.Item(132, 13).Value = -1970
$ws.Cells.Item(132, 14).Value = -14033
$ws.Cells.Item(134, 8).Value = 3593.88
$ws.Cells.Item(134, 9).Value = 3750.35
$ws.Cells.Item(134, 11).Value = 11251.05
$ws.Cells.Item(134, 13).Value = -8716.049999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 998691.4399999999
$ws.Cells.Item(102, 9).Value = 1884560
$ws.Cells.Item(102, 10).Value = 2089.25
$ws.Cells.Item(102, 11).Value = 1884560
$ws.Cells.Item(102, 12).Value = 2089.25
$ws.Cells.Item(102, 13).Value = -1882938
$ws.Cells.Item(102, 14).Value = -5333.25
$ws.Cells.Item(122, 8).Value = 796042.2
$ws.Cells.Item(122, 9).Value = 1264408.5
$ws.Cells.Item(122, 10).Value = 3422.1538
$ws.Cells.Item(122, 11).Value = 3793225.5
$ws.Cells.Item(122, 12).Value = 10266.4614
$ws.Cells.Item(122, 13).Value = -3790775.5
$ws.Cells.Item(122, 14).Value = -15166.4614
$ws.Cells.Item(126, 8).Value = 4331.896
$ws.Cells.Item(126, 9).Value = 6276.905
$ws.Cells.Item(126, 10).Value = 2819.111
$ws.Cells.Item(126, 11).Value = 18830.715
$ws.Cells.Item(126, 12).Value = 8457.332999999999
$ws.Cells.Item(126, 13).Value = -16360.715
$ws.Cells.Item(126, 14).Value = -13397.333
$ws.Cells.Item(132, 8).Value = 2756.85
$ws.Cells.Item(132, 9).Value = 2289.5
$ws.Cells.Item(132, 10).Value = 3068.4167
$ws.Cells.Item(132, 11).Value = 6868.5
$ws.Cells.Item(132, 12).Value = 9205.250100000001
$ws.Cells.Item(132, 13).Value = -4338.5
$ws.Cells.Item(132, 14).Value = -14265.2501
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 34913.867
$ws.Cells.Item(7, 9).Value = 51645.3
$ws.Cells.Item(7, 10).Value = 1451
$ws.Cells.Item(7, 11).Value = 51645.3
$ws.Cells.Item(7, 12).Value = 1451
$ws.Cells.Item(7, 13).Value = -51533.3
$ws.Cells.Item(7, 14).Value = -1675
$ws.Cells.Item(40, 8).Value = 55558464
$ws.Cells.Item(40, 9).Value = 100001920
$ws.Cells.Item(40, 10).Value = 4143.125
$ws.Cells.Item(40, 11).Value = 100001920
$ws.Cells.Item(40, 12).Value = 4143.125
$ws.Cells.Item(40, 13).Value = -100001784
$ws.Cells.Item(40, 14).Value = -4415.125
$ws.Cells.Item(74, 8).Value = 19666.666
$ws.Cells.Item(74, 10).Value = 19666.666
$ws.Cells.Item(74, 12).Value = 19666.666
$ws.Cells.Item(74, 14).Value = -21662.666
$ws.Cells.Item(77, 8).Value = 19666.666
$ws.Cells.Item(77, 10).Value = 19666.666
$ws.Cells.Item(77, 12).Value = 58999.99800000001
$ws.Cells.Item(77, 14).Value = -68983.99800000001
$ws.Cells.Item(126, 8).Value = 34913.867
$ws.Cells.Item(126, 9).Value = 51645.3
$ws.Cells.Item(126, 10).Value = 1451
$ws.Cells.Item(126, 11).Value = 154935.9
$ws.Cells.Item(126, 12).Value = 4353
$ws.Cells.Item(126, 13).Value = -152465.9
$ws.Cells.Item(126, 14).Value = -9293
$ws.Cells.Item(132, 8).Value = 12351318
$ws.Cells.Item(132, 9).Value = 18525106
$ws.Cells.Item(132, 10).Value = 3742.889
$ws.Cells.Item(132, 11).Value = 55575318
$ws.Cells.Item(132, 12).Value = 11228.667
$ws.Cells.Item(132, 13).Value = -55572788
$ws.Cells.Item(132, 14).Value = -16288.667
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 66667360
$ws.Cells.Item(107, 9).Value = 125000744
$ws.Cells.Item(107, 10).Value = 635.4286
$ws.Cells.Item(107, 11).Value = 375002232
$ws.Cells.Item(107, 12).Value = 1906.2858
$ws.Cells.Item(107, 13).Value = -375000312
$ws.Cells.Item(107, 14).Value = -5746.2858
$ws.Cells.Item(122, 8).Value = 1070
$ws.Cells.Item(122, 9).Value = 860
$ws.Cells.Item(122, 10).Value = 1700
$ws.Cells.Item(122, 11).Value = 2580
$ws.Cells.Item(122, 12).Value = 5100
$ws.Cells.Item(122, 13).Value = -130
$ws.Cells.Item(122, 14).Value = -10000
$ws.Cells.Item(132, 8).Value = 1127.6852
$ws.Cells.Item(132, 9).Value = 785.61536
$ws.Cells.Item(132, 11).Value = 2356.84608
$ws.Cells.Item(132, 13).Value = 173.1539199999997
$ws.Cells.Item(136, 8).Value = 2012.4073
$ws.Cells.Item(136, 9).Value = 1958.7097
$ws.Cells.Item(136, 10).Value = 2084.7827
$ws.Cells.Item(136, 11).Value = 5876.1291
$ws.Cells.Item(136, 12).Value = 6254.348100000001
$ws.Cells.Item(136, 13).Value = -3326.1291
$ws.Cells.Item(136, 14).Value = -11354.3481
